$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CLS")

# Step 1: Insert two new columns before column D (shifts D:K -> F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Step 2: Copy number formatting from the (old D:E, now F:G) columns into the
# newly inserted D:E columns so date/number styles match the rest of the table.
# Restricted to the exact row blocks that contain populated D:M cells, so we
# do not manufacture styled-but-empty cells on header/spacer rows.
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: Write the new Dec-2018 / Sep-2018 quarter data (columns D, E) and a
# handful of restated historical values (within columns F-J) for CLS quarterly financials.
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 4).Value = 1727000
$ws.Cells.Item(8, 5).Value = 1711300
$ws.Cells.Item(8, 8).Value = 1570200
$ws.Cells.Item(8, 9).Value = 1532800
$ws.Cells.Item(9, 4).Value = 1607000
$ws.Cells.Item(9, 5).Value = 1599100
$ws.Cells.Item(9, 6).Value = 1588800
$ws.Cells.Item(9, 8).Value = 2920100
$ws.Cells.Item(9, 9).Value = 2851300
$ws.Cells.Item(9, 10).Value = 2898000
$ws.Cells.Item(10, 4).Value = 120000
$ws.Cells.Item(10, 5).Value = 112200
$ws.Cells.Item(10, 6).Value = 106400
$ws.Cells.Item(10, 8).Value = -1349900
$ws.Cells.Item(10, 9).Value = -1318500
$ws.Cells.Item(10, 10).Value = -1340400
$ws.Cells.Item(12, 4).Value = 8100
$ws.Cells.Item(12, 5).Value = 7900
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = 16900
$ws.Cells.Item(14, 5).Value = 17900
$ws.Cells.Item(14, 6).Value = 16100
$ws.Cells.Item(14, 7).Value = 8600
$ws.Cells.Item(14, 8).Value = 17500
$ws.Cells.Item(14, 9).Value = 3900
$ws.Cells.Item(14, 10).Value = 7500
$ws.Cells.Item(15, 4).Value = 6100
$ws.Cells.Item(15, 5).Value = 3600
$ws.Cells.Item(17, 4).Value = 1697700
$ws.Cells.Item(17, 5).Value = 1682800
$ws.Cells.Item(17, 8).Value = 1546300
$ws.Cells.Item(17, 9).Value = 1488200
$ws.Cells.Item(18, 4).Value = 29300
$ws.Cells.Item(18, 5).Value = 28500
$ws.Cells.Item(18, 8).Value = 23900
$ws.Cells.Item(18, 9).Value = 44600
$ws.Cells.Item(20, 4).Value = -9200
$ws.Cells.Item(20, 5).Value = -7000
$ws.Cells.Item(21, 4).Value = 45100
$ws.Cells.Item(21, 5).Value = 42700
$ws.Cells.Item(21, 8).Value = 40900
$ws.Cells.Item(21, 9).Value = 61700
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(23, 4).Value = 20100
$ws.Cells.Item(23, 5).Value = 21500
$ws.Cells.Item(23, 8).Value = 21300
$ws.Cells.Item(23, 9).Value = 42300
$ws.Cells.Item(24, 4).Value = -40000
$ws.Cells.Item(24, 5).Value = 12900
$ws.Cells.Item(24, 9).Value = 7500
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 60100
$ws.Cells.Item(26, 5).Value = 8600
$ws.Cells.Item(26, 8).Value = 15600
$ws.Cells.Item(26, 9).Value = 34800
$ws.Cells.Item(27, 4).Value = 60100
$ws.Cells.Item(27, 5).Value = 8600
$ws.Cells.Item(27, 8).Value = 15600
$ws.Cells.Item(27, 9).Value = 34800
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = "NA"
$ws.Cells.Item(29, 5).Value = "NA"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = 9200
$ws.Cells.Item(32, 5).Value = 7000
$ws.Cells.Item(33, 4).Value = 60100
$ws.Cells.Item(33, 5).Value = 8600
$ws.Cells.Item(33, 8).Value = 13600
$ws.Cells.Item(33, 9).Value = 34800
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = 60100
$ws.Cells.Item(35, 5).Value = 8600
$ws.Cells.Item(35, 8).Value = 13600
$ws.Cells.Item(35, 9).Value = 34800
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(41, 4).Value = 422000
$ws.Cells.Item(41, 5).Value = 457700
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(43, 4).Value = 1211600
$ws.Cells.Item(43, 5).Value = 1135700
$ws.Cells.Item(44, 4).Value = 1089900
$ws.Cells.Item(44, 5).Value = 1062400
$ws.Cells.Item(45, 4).Value = 100000
$ws.Cells.Item(45, 5).Value = 95500
$ws.Cells.Item(46, 4).Value = 2823500
$ws.Cells.Item(46, 5).Value = 2751300
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(48, 4).Value = 365300
$ws.Cells.Item(48, 5).Value = 344300
$ws.Cells.Item(49, 4).Value = 482000
$ws.Cells.Item(49, 5).Value = 152800
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 66900
$ws.Cells.Item(52, 5).Value = 67700
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 3737700
$ws.Cells.Item(54, 5).Value = 3316100
$ws.Cells.Item(57, 4).Value = 1126700
$ws.Cells.Item(57, 5).Value = 1147800
$ws.Cells.Item(58, 4).Value = 107700
$ws.Cells.Item(58, 5).Value = 61000
$ws.Cells.Item(59, 4).Value = 385900
$ws.Cells.Item(59, 5).Value = 325500
$ws.Cells.Item(60, 4).Value = 1620300
$ws.Cells.Item(60, 5).Value = 1534300
$ws.Cells.Item(61, 4).Value = 650200
$ws.Cells.Item(61, 5).Value = 350600
$ws.Cells.Item(62, 4).Value = 134900
$ws.Cells.Item(62, 5).Value = 142600
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 2405400
$ws.Cells.Item(66, 5).Value = 2027500
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = -1481700
$ws.Cells.Item(72, 5).Value = -1550200
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 1332300
$ws.Cells.Item(76, 5).Value = 1288600
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 4).Value = 60100
$ws.Cells.Item(81, 5).Value = 8600
$ws.Cells.Item(81, 8).Value = 13600
$ws.Cells.Item(81, 9).Value = 34800
$ws.Cells.Item(83, 4).Value = 25000
$ws.Cells.Item(83, 5).Value = 21200
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = -1900
$ws.Cells.Item(89, 5).Value = 55300
$ws.Cells.Item(91, 4).Value = -18800
$ws.Cells.Item(91, 5).Value = -21100
$ws.Cells.Item(91, 9).Value = -32200
$ws.Cells.Item(91, 10).Value = -24100
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = -344200
$ws.Cells.Item(94, 5).Value = -20900
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(96, 5).Value = 0
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = 310400
$ws.Cells.Item(100, 5).Value = 21900
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(101, 5).Value = 0
$ws.Cells.Item(102, 4).Value = -35700
$ws.Cells.Item(102, 5).Value = 56300

$ws.Range("A5:M102").Columns.AutoFit() | Out-Null

Write-Host "Applied CLS quarterly financial updates"
